# Auto-generated PowerShell COM-interop script
# Commit: feat: add 2022-Q1 data
#
# Plan (mirrors the OOXML diff):
#  1. The existing "总计" sheet (sheetId 6 / sheet6.xml) is renamed to
#     "2022-Q1" and repopulated with the fund-holding detail table.
#  2. A brand-new "总计" sheet is inserted right after it (it naturally
#     receives the next sheetId, 7 / sheet7.xml) with the summary table,
#     now including a "2022-Q1" row.
$wb = $excel.ActiveWorkbook

$wsQ1 = $wb.Worksheets.Item(6)
$wsQ1.Name = "2022-Q1"
$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

# Canonical "header / index" style (bold, centered, thin border) already
# lives on $wsQ1 (it is what the old "总计" sheet used) - reuse it via
# copy/paste-special (formats only) instead of rebuilding it from scratch,
# so every new cell matches the style used throughout the workbook.
$styleSrc = $wsQ1.Range("A2")

# --- "2022-Q1" sheet: fund holding detail table (A1:H45) -----------------
# B1:D1 already carry the bordered header style inherited from the old
# "总计" sheet; E1:H1 are brand new cells so copy the style onto them first.
$styleSrc.Copy()
$wsQ1.Range("E1:H1").PasteSpecial(-4122)
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Column A (row index) keeps its existing style for rows 2-6; copy it down
# for the newly added rows 7-45.
$styleSrc.Copy()
$wsQ1.Range("A7:A45").PasteSpecial(-4122)

# Columns B (fund code) and D:G (decimal-looking scale/position/pct/value)
# are stored as literal TEXT in this workbook (leading zeros in fund codes,
# trailing zeros in the decimals) - force Text format before assigning so
# Excel does not coerce them into numbers.
$wsQ1.Range("B2:B45").NumberFormat = "@"
$wsQ1.Range("D2:G45").NumberFormat = "@"

$q1Data = @(
  @(0, "011174", "中庚价值品质一年持有期混合", "63.01", "93.47", "6.57", "4.1398", 6),
  @(1, "008954", "安信价值回报三年持有期混合A", "33.45", "90.20", "6.99", "2.3382", 5),
  @(2, "006551", "中庚价值领航混合", "36.49", "94.18", "4.56", "1.6639", 7),
  @(3, "000934", "国富大中华精选混合QDII", "25.71", "83.59", "3.24", "0.8330", 4),
  @(4, "006370", "国富大中华精选混合QDII美元", "25.71", "83.59", "3.24", "0.8330", 4),
  @(5, "008809", "安信民稳增长混合A", "15.09", "44.77", "4.98", "0.7515", 1),
  @(6, "011128", "华安精致生活混合A", "33.30", "85.22", "2.10", "0.6993", 8),
  @(7, "009983", "永赢港股通品质生活慧选混合", "9.75", "87.97", "6.63", "0.6464", 2),
  @(8, "012256", "安信丰穗一年持有混合A", "26.49", "20.42", "1.97", "0.5219", 1),
  @(9, "009618", "交银施罗德启汇混合", "21.00", "71.65", "2.23", "0.4683", 10),
  @(10, "001583", "安信新常态沪港深精选股票", "8.85", "89.85", "5.17", "0.4575", 5),
  @(11, "009849", "安信稳健聚申一年持有期混合", "12.33", "30.81", "3.47", "0.4279", 1),
  @(12, "167508", "安信价值发现两年定期开放混合（LOF）", "5.43", "89.42", "7.64", "0.4149", 1),
  @(13, "008810", "安信民稳增长混合C", "6.76", "44.77", "4.98", "0.3366", 1),
  @(14, "011905", "安信价值启航混合型证券投资基金A", "4.48", "93.14", "6.38", "0.2858", 3),
  @(15, "013123", "汇添富精选核心优势一年持有混合A", "6.15", "66.61", "4.45", "0.2737", 2),
  @(16, "011315", "永赢港股通优质成长一年持有期混合型证券投资基金", "4.07", "90.36", "6.54", "0.2662", 1),
  @(17, "011913", "华夏永泓一年持有混合A", "24.65", "37.51", "1.04", "0.2564", 7),
  @(18, "012073", "华安均衡优选混合A", "8.33", "89.23", "3.04", "0.2532", 9),
  @(19, "519779", "交银施罗德沪港深价值精选灵活配置混合", "5.13", "84.44", "4.45", "0.2283", 8),
  @(20, "001581", "华安沪港深通精选灵活配置混合", "4.92", "92.91", "3.56", "0.1752", 6),
  @(21, "011129", "华安精致生活混合C", "7.57", "85.22", "2.10", "0.1590", 8),
  @(22, "008477", "安信价值驱动三年持有期混合", "2.64", "91.75", "6.00", "0.1584", 3),
  @(23, "010204", "中银港股通优势成长股票", "3.19", "83.00", "4.77", "0.1522", 6),
  @(24, "011914", "华夏永泓一年持有混合C", "12.39", "37.51", "1.04", "0.1289", 7),
  @(25, "005259", "建信龙头企业股票", "1.04", "83.45", "6.76", "0.0703", 7),
  @(26, "011969", "建信港股通精选混合A", "1.01", "57.64", "6.79", "0.0686", 2),
  @(27, "011355", "华泰柏瑞港股通时代机遇混合型证券投资基金A", "1.13", "90.93", "5.87", "0.0663", 8),
  @(28, "004249", "安信中国制造2025沪港深灵活配置混合", "0.58", "89.89", "8.32", "0.0483", 2),
  @(29, "012257", "安信丰穗一年持有混合C", "2.41", "20.42", "1.97", "0.0475", 1),
  @(30, "003413", "华泰柏瑞新经济沪港深灵活配置混合", "0.54", "92.57", "7.98", "0.0431", 5),
  @(31, "007674", "工银瑞信产业升级股票A", "0.69", "94.04", "6.04", "0.0417", 6),
  @(32, "004393", "安信合作创新主题沪港深灵活配置混合", "0.49", "89.26", "7.90", "0.0387", 3),
  @(33, "005142", "中融沪港深大消费主题灵活配置混合A", "0.49", "88.98", "5.62", "0.0275", 2),
  @(34, "011906", "安信价值启航混合型证券投资基金C", "0.42", "93.14", "6.38", "0.0268", 3),
  @(35, "011356", "华泰柏瑞港股通时代机遇混合型证券投资基金C", "0.40", "90.93", "5.87", "0.0235", 8),
  @(36, "011970", "建信港股通精选混合C", "0.33", "57.64", "6.79", "0.0224", 2),
  @(37, "002214", "中海沪港深价值优选灵活配置混合", "0.65", "84.71", "3.11", "0.0202", 8),
  @(38, "007675", "工银瑞信产业升级股票C", "0.33", "94.04", "6.04", "0.0199", 6),
  @(39, "005143", "中融沪港深大消费主题灵活配置混合C", "0.33", "88.98", "5.62", "0.0185", 2),
  @(40, "013124", "汇添富精选核心优势一年持有混合C", "0.30", "66.61", "4.45", "0.0134", 2),
  @(41, "012074", "华安均衡优选混合C", "0.25", "89.23", "3.04", "0.0076", 9),
  @(42, "040021", "华安大中华升级股票(QDII)", "0.26", "87.37", "2.29", "0.0060", 10),
  @(43, "010667", "安信价值回报三年持有期混合C", "0.02", "90.20", "6.99", "0.0014", 5)
)

$r = 2
foreach ($row in $q1Data) {
    $wsQ1.Cells.Item($r, 1).Value = $row[0]
    $wsQ1.Cells.Item($r, 2).Value = $row[1]
    $wsQ1.Cells.Item($r, 3).Value = $row[2]
    $wsQ1.Cells.Item($r, 4).Value = $row[3]
    $wsQ1.Cells.Item($r, 5).Value = $row[4]
    $wsQ1.Cells.Item($r, 6).Value = $row[5]
    $wsQ1.Cells.Item($r, 7).Value = $row[6]
    $wsQ1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- "总计" sheet: quarterly summary table (A1:D7), now with 2022-Q1 -----
$styleSrc.Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)
$wsTotal.Range("A2:A7").PasteSpecial(-4122)
$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
  @(0, "2022-Q1", 44, 17.48),
  @(1, "2021-Q4", 26, 12.09),
  @(2, "2021-Q3", 21, 10.45),
  @(3, "2021-Q2", 12, 6.98),
  @(4, "2021-Q1", 18, 5.78),
  @(5, "2020-Q4", 2, 0.02)
)

$r = 2
foreach ($row in $totalData) {
    $wsTotal.Cells.Item($r, 1).Value = $row[0]
    $wsTotal.Cells.Item($r, 2).Value = $row[1]
    $wsTotal.Cells.Item($r, 3).Value = $row[2]
    $wsTotal.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$wsTotal.Range("A1").Select()
